# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-03 14:16:11
#
# Applies the diffs to "Session Analysis Results":
#  - reorders/updates several "Recorded By" email lists
#  - updates recorded/missing session counts + coverage/attendance percentages
#  - flips three sessions (rows 19, 41, 76) from "Not Recorded" (pink) to
#    "Recorded" (green), filling in their Recorded-By / Students / Status cells
#  - updates a handful of per-group stat numbers + percentages (cols O/P/R/S)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel re-interpreting
# it as a number/date/percentage (e.g. "17.6%" or "43/216"), and without
# leaving a stray "quote prefix" / number-format behind - we borrow the
# destination's original formatting from a still-untouched cell that already
# carries the same "s=5" style (L4, a plain numeric stat cell nothing in this
# edit ever touches), then paste formats (not values) back on top.
# ---------------------------------------------------------------------------
$script:FormatDonor = "L4"

function Set-TextValue {
    param($targetAddr, [string]$text)

    $ws.Range($targetAddr).Value = "'" + $text
    $ws.Range($script:FormatDonor).Copy()
    $ws.Range($targetAddr).PasteSpecial(-4122)   # xlPasteFormats
}

# ---------------------------------------------------------------------------
# Recorded-By (column G) email list reorders / additions
# ---------------------------------------------------------------------------
$ws.Range("G2").Value  = "shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G24").Value = "shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"

$ws.Range("G18").Value = "yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G40").Value = "yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"

$ws.Range("G54").Value  = "merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G98").Value  = "merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G120").Value = "merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"

$ws.Range("G58").Value = "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"
$ws.Range("G80").Value = "afaf.abdallah@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg"

$ws.Range("G62").Value = "wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G84").Value = "wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"

$ws.Range("G96").Value  = "Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"
$ws.Range("G118").Value = "Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg"

$ws.Range("G106").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G128").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"

$ws.Range("G134").Value = "hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

$ws.Range("G150").Value = "Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G172").Value = "Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"

$ws.Range("G156").Value = "alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# ---------------------------------------------------------------------------
# Students (col H) count updates that accompany the above Recorded-By edits
# ---------------------------------------------------------------------------
$ws.Range("H54").Value  = "71/220"
$ws.Range("H106").Value = "61/154"
$ws.Range("H172").Value = "50/226"

# ---------------------------------------------------------------------------
# Class Statistics summary block (K3:L10)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 11
Set-TextValue "L9" "17.6%"
Set-TextValue "L10" "29.8%"

# ---------------------------------------------------------------------------
# Per-group statistics block (cols K:S) - numeric + percentage tweaks
# ---------------------------------------------------------------------------
$ws.Range("O15").Value = 4
$ws.Range("P15").Value = 2
Set-TextValue "R15" "18.2%"
Set-TextValue "S15" "32.2%"

$ws.Range("O16").Value = 4
$ws.Range("P16").Value = 2
Set-TextValue "R16" "18.2%"
Set-TextValue "S16" "28.9%"

Set-TextValue "S17" "40.8%"

$ws.Range("O18").Value = 4
$ws.Range("P18").Value = 1
Set-TextValue "R18" "18.2%"
Set-TextValue "S18" "25.8%"

Set-TextValue "S19" "29.2%"
Set-TextValue "S22" "16.8%"

# ---------------------------------------------------------------------------
# Rows 19 / 41 / 76: sessions flip from "Not Recorded" (pink) to "Recorded"
# (green). Copy the formatting of an already-"Recorded" row in the same
# block onto columns A:I, then fill in the real Recorded-By / Students /
# Status values (the other A:F fields keep their existing values).
# ---------------------------------------------------------------------------
$ws.Range("A18:I18").Copy()
$ws.Range("A19:I19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G19").Value = "Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("H19").Value = "43/216"
$ws.Range("I19").Value = "Recorded"

$ws.Range("A40:I40").Copy()
$ws.Range("A41:I41").PasteSpecial(-4122)
$ws.Range("G41").Value = "Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("H41").Value = "56/217"
$ws.Range("I41").Value = "Recorded"

$ws.Range("A54:I54").Copy()
$ws.Range("A76:I76").PasteSpecial(-4122)   # copies formatting only, values untouched
$ws.Range("G76").Value = "merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("H76").Value = "50/225"
$ws.Range("I76").Value = "Recorded"
